$d = $word.ActiveDocument

# 1) Update the Revit SDK year references from 2018 to 2024.
#    (Appears twice: once in the "Extensions 2018\Framework\Foundation"
#    path and once in the "Addins\2018" path.)
$rng = $d.Content
$rng.Find.Execute("2018", $false, $false, $false, $false, $false, $true, 1, $false, "2024", 2) | Out-Null

# 2) Remove the leftover "_GoBack" bookmark paragraph at the end of the
#    document (Word drops this automatically-managed bookmark once the
#    document is edited and resaved).
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
} catch {
}
